# Weekly fruit/vegetable price update: insert 3 new data rows for the
# latest reporting week (2023-12-07, serial 45267) at the top of the
# "Cereza" (cherry) records block, pushing the existing rows down.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows right before the existing data block (row 328),
# shifting rows 328:350 down to 331:353.
$ws.Rows("328:330").Insert()

# New row 328: Early Burlat, Primera
$ws.Cells.Item(328, 1).Value = 5
$ws.Cells.Item(328, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(328, 3).Value = "Maule"
$ws.Cells.Item(328, 4).Value = 45267
$ws.Cells.Item(328, 5).Value = 7
$ws.Cells.Item(328, 6).Value = "Fruta"
$ws.Cells.Item(328, 7).Value = 100103
$ws.Cells.Item(328, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(328, 9).Value = 100103001
$ws.Cells.Item(328, 10).Value = "Cereza"
$ws.Cells.Item(328, 11).Value = "Early Burlat"
$ws.Cells.Item(328, 12).Value = "Primera"
$ws.Cells.Item(328, 13).Value = 150
$ws.Cells.Item(328, 14).Value = 7000
$ws.Cells.Item(328, 15).Value = 7000
$ws.Cells.Item(328, 16).Value = 7000
$ws.Cells.Item(328, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(328, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(328, 19).Value = 700
$ws.Cells.Item(328, 20).Value = 10

# New row 329: Rainier, Primera
$ws.Cells.Item(329, 1).Value = 5
$ws.Cells.Item(329, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(329, 3).Value = "Maule"
$ws.Cells.Item(329, 4).Value = 45267
$ws.Cells.Item(329, 5).Value = 7
$ws.Cells.Item(329, 6).Value = "Fruta"
$ws.Cells.Item(329, 7).Value = 100103
$ws.Cells.Item(329, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(329, 9).Value = 100103001
$ws.Cells.Item(329, 10).Value = "Cereza"
$ws.Cells.Item(329, 11).Value = "Rainier"
$ws.Cells.Item(329, 12).Value = "Primera"
$ws.Cells.Item(329, 13).Value = 100
$ws.Cells.Item(329, 14).Value = 10000
$ws.Cells.Item(329, 15).Value = 10000
$ws.Cells.Item(329, 16).Value = 10000
$ws.Cells.Item(329, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(329, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(329, 19).Value = 1000
$ws.Cells.Item(329, 20).Value = 10

# New row 330: Santina, Primera
$ws.Cells.Item(330, 1).Value = 5
$ws.Cells.Item(330, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(330, 3).Value = "Maule"
$ws.Cells.Item(330, 4).Value = 45267
$ws.Cells.Item(330, 5).Value = 7
$ws.Cells.Item(330, 6).Value = "Fruta"
$ws.Cells.Item(330, 7).Value = 100103
$ws.Cells.Item(330, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(330, 9).Value = 100103001
$ws.Cells.Item(330, 10).Value = "Cereza"
$ws.Cells.Item(330, 11).Value = "Santina"
$ws.Cells.Item(330, 12).Value = "Primera"
$ws.Cells.Item(330, 13).Value = 120
$ws.Cells.Item(330, 14).Value = 7000
$ws.Cells.Item(330, 15).Value = 7000
$ws.Cells.Item(330, 16).Value = 7000
$ws.Cells.Item(330, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(330, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(330, 19).Value = 700
$ws.Cells.Item(330, 20).Value = 10
